$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: fill column G (time values) for rows 528-547 first
# This establishes shared-string order for indices 504-522
$ws.Cells.Item(528,7).Value = '01:41:34'
$ws.Cells.Item(529,7).Value = '01:40:12'
$ws.Cells.Item(530,7).Value = '00:18:53'
$ws.Cells.Item(531,7).Value = '00:11:49'
$ws.Cells.Item(532,7).Value = '01:22:09'
$ws.Cells.Item(533,7).Value = '00:30:08'
$ws.Cells.Item(534,7).Value = '01:29:04'
$ws.Cells.Item(535,7).Value = '01:41:10'
$ws.Cells.Item(536,7).Value = '01:41:18'
$ws.Cells.Item(537,7).Value = '01:41:18'
$ws.Cells.Item(538,7).Value = '00:19:01'
$ws.Cells.Item(539,7).Value = '00:21:11'
$ws.Cells.Item(540,7).Value = '01:11:58'
$ws.Cells.Item(541,7).Value = '01:21:36'
$ws.Cells.Item(542,7).Value = '01:36:17'
$ws.Cells.Item(543,7).Value = '01:33:29'
$ws.Cells.Item(544,7).Value = '01:25:21'
$ws.Cells.Item(545,7).Value = '01:34:28'
$ws.Cells.Item(546,7).Value = '01:32:55'
$ws.Cells.Item(547,7).Value = '01:29:17'

# Step 2: fill column A (match name) for rows 528-541
# This establishes shared-string order for index 523
$ws.Cells.Item(528,1).Value = 'N3 J4 VS Bourgoin'
$ws.Cells.Item(529,1).Value = 'N3 J4 VS Bourgoin'
$ws.Cells.Item(530,1).Value = 'N3 J4 VS Bourgoin'
$ws.Cells.Item(531,1).Value = 'N3 J4 VS Bourgoin'
$ws.Cells.Item(532,1).Value = 'N3 J4 VS Bourgoin'
$ws.Cells.Item(533,1).Value = 'N3 J4 VS Bourgoin'
$ws.Cells.Item(534,1).Value = 'N3 J4 VS Bourgoin'
$ws.Cells.Item(535,1).Value = 'N3 J4 VS Bourgoin'
$ws.Cells.Item(536,1).Value = 'N3 J4 VS Bourgoin'
$ws.Cells.Item(537,1).Value = 'N3 J4 VS Bourgoin'
$ws.Cells.Item(538,1).Value = 'N3 J4 VS Bourgoin'
$ws.Cells.Item(539,1).Value = 'N3 J4 VS Bourgoin'
$ws.Cells.Item(540,1).Value = 'N3 J4 VS Bourgoin'
$ws.Cells.Item(541,1).Value = 'N3 J4 VS Bourgoin'

# Step 3: fill column G for rows 548-552 (establishes indices 524-528)
$ws.Cells.Item(548,7).Value = '01:26:06'
$ws.Cells.Item(549,7).Value = '01:24:28'
$ws.Cells.Item(550,7).Value = '01:20:38'
$ws.Cells.Item(551,7).Value = '01:19:53'
$ws.Cells.Item(552,7).Value = '01:23:14'

# Step 4: fill column A for rows 542-552 (reuses existing string 'Entrainement')
$ws.Cells.Item(542,1).Value = 'Entrainement'
$ws.Cells.Item(543,1).Value = 'Entrainement'
$ws.Cells.Item(544,1).Value = 'Entrainement'
$ws.Cells.Item(545,1).Value = 'Entrainement'
$ws.Cells.Item(546,1).Value = 'Entrainement'
$ws.Cells.Item(547,1).Value = 'Entrainement'
$ws.Cells.Item(548,1).Value = 'Entrainement'
$ws.Cells.Item(549,1).Value = 'Entrainement'
$ws.Cells.Item(550,1).Value = 'Entrainement'
$ws.Cells.Item(551,1).Value = 'Entrainement'
$ws.Cells.Item(552,1).Value = 'Entrainement'

# Step 5: fill remaining text columns B, C, D, E, F and numeric columns H-V for all new rows
$ws.Cells.Item(528,2).Value = 45920
$ws.Cells.Item(528,3).Value = 'Global'
$ws.Cells.Item(528,4).Value = 'M'
$ws.Cells.Item(528,5).Value = 'Naim Dhib'
$ws.Cells.Item(528,6).Value = 'center midfield'
$ws.Cells.Item(528,8).Value = 10.19
$ws.Cells.Item(528,9).Value = 1.5
$ws.Cells.Item(528,10).Value = 8.67
$ws.Cells.Item(528,11).Value = 1.08
$ws.Cells.Item(528,12).Value = 0.31
$ws.Cells.Item(528,13).Value = 0.13
$ws.Cells.Item(528,14).Value = 0
$ws.Cells.Item(528,15).Value = 10
$ws.Cells.Item(528,16).Value = 5.98
$ws.Cells.Item(528,17).Value = 28.78
$ws.Cells.Item(528,18).Value = 4.99
$ws.Cells.Item(528,19).Value = 32
$ws.Cells.Item(528,20).Value = 11
$ws.Cells.Item(528,21).Value = 42
$ws.Cells.Item(528,22).Value = 15

$ws.Cells.Item(529,2).Value = 45920
$ws.Cells.Item(529,3).Value = 'Global'
$ws.Cells.Item(529,4).Value = 'M'
$ws.Cells.Item(529,5).Value = 'Yoan Zouma'
$ws.Cells.Item(529,6).Value = 'center back'
$ws.Cells.Item(529,8).Value = 8.71
$ws.Cells.Item(529,9).Value = 0.96
$ws.Cells.Item(529,10).Value = 7.74
$ws.Cells.Item(529,11).Value = 0.68
$ws.Cells.Item(529,12).Value = 0.2
$ws.Cells.Item(529,13).Value = 0.09
$ws.Cells.Item(529,14).Value = 0
$ws.Cells.Item(529,15).Value = 6
$ws.Cells.Item(529,16).Value = 5.19
$ws.Cells.Item(529,17).Value = 29.17
$ws.Cells.Item(529,18).Value = 4.41
$ws.Cells.Item(529,19).Value = 26
$ws.Cells.Item(529,20).Value = 7
$ws.Cells.Item(529,21).Value = 26
$ws.Cells.Item(529,22).Value = 11

$ws.Cells.Item(530,2).Value = 45920
$ws.Cells.Item(530,3).Value = 'Global'
$ws.Cells.Item(530,4).Value = 'M'
$ws.Cells.Item(530,5).Value = 'Amir Etien'
$ws.Cells.Item(530,6).Value = 'right forward'
$ws.Cells.Item(530,8).Value = 1.77
$ws.Cells.Item(530,9).Value = 0.36
$ws.Cells.Item(530,10).Value = 1.4
$ws.Cells.Item(530,11).Value = 0.18
$ws.Cells.Item(530,12).Value = 0.08
$ws.Cells.Item(530,13).Value = 0.09
$ws.Cells.Item(530,14).Value = 0.02
$ws.Cells.Item(530,15).Value = 4
$ws.Cells.Item(530,16).Value = 5.73
$ws.Cells.Item(530,17).Value = 31.99
$ws.Cells.Item(530,18).Value = 4.15
$ws.Cells.Item(530,19).Value = 8
$ws.Cells.Item(530,20).Value = 1
$ws.Cells.Item(530,21).Value = 8
$ws.Cells.Item(530,22).Value = 3

$ws.Cells.Item(531,2).Value = 45920
$ws.Cells.Item(531,3).Value = 'Global'
$ws.Cells.Item(531,4).Value = 'M'
$ws.Cells.Item(531,5).Value = 'Emmanuel Valey'
$ws.Cells.Item(531,6).Value = 'left forward'
$ws.Cells.Item(531,8).Value = 1.42
$ws.Cells.Item(531,9).Value = 0.36
$ws.Cells.Item(531,10).Value = 1.06
$ws.Cells.Item(531,11).Value = 0.23
$ws.Cells.Item(531,12).Value = 0.1
$ws.Cells.Item(531,13).Value = 0.04
$ws.Cells.Item(531,14).Value = 0
$ws.Cells.Item(531,15).Value = 5
$ws.Cells.Item(531,16).Value = 7.19
$ws.Cells.Item(531,17).Value = 29.26
$ws.Cells.Item(531,18).Value = 3.82
$ws.Cells.Item(531,19).Value = 11
$ws.Cells.Item(531,20).Value = 0
$ws.Cells.Item(531,21).Value = 7
$ws.Cells.Item(531,22).Value = 0

$ws.Cells.Item(532,2).Value = 45920
$ws.Cells.Item(532,3).Value = 'Global'
$ws.Cells.Item(532,4).Value = 'M'
$ws.Cells.Item(532,5).Value = 'Karahali Souaré'
$ws.Cells.Item(532,6).Value = 'right forward'
$ws.Cells.Item(532,8).Value = 8.41
$ws.Cells.Item(532,9).Value = 1.57
$ws.Cells.Item(532,10).Value = 6.82
$ws.Cells.Item(532,11).Value = 0.97
$ws.Cells.Item(532,12).Value = 0.49
$ws.Cells.Item(532,13).Value = 0.13
$ws.Cells.Item(532,14).Value = 0
$ws.Cells.Item(532,15).Value = 10
$ws.Cells.Item(532,16).Value = 6.1
$ws.Cells.Item(532,17).Value = 28.66
$ws.Cells.Item(532,18).Value = 5.07
$ws.Cells.Item(532,19).Value = 51
$ws.Cells.Item(532,20).Value = 15
$ws.Cells.Item(532,21).Value = 43
$ws.Cells.Item(532,22).Value = 16

$ws.Cells.Item(533,2).Value = 45920
$ws.Cells.Item(533,3).Value = 'Global'
$ws.Cells.Item(533,4).Value = 'M'
$ws.Cells.Item(533,5).Value = 'Malik Boussaid'
$ws.Cells.Item(533,6).Value = 'right back'
$ws.Cells.Item(533,8).Value = 3.56
$ws.Cells.Item(533,9).Value = 0.68
$ws.Cells.Item(533,10).Value = 2.87
$ws.Cells.Item(533,11).Value = 0.5
$ws.Cells.Item(533,12).Value = 0.17
$ws.Cells.Item(533,13).Value = 0.02
$ws.Cells.Item(533,14).Value = 0
$ws.Cells.Item(533,15).Value = 2
$ws.Cells.Item(533,16).Value = 7.04
$ws.Cells.Item(533,17).Value = 27.6
$ws.Cells.Item(533,18).Value = 3.97
$ws.Cells.Item(533,19).Value = 13
$ws.Cells.Item(533,20).Value = 0
$ws.Cells.Item(533,21).Value = 22
$ws.Cells.Item(533,22).Value = 8

$ws.Cells.Item(534,2).Value = 45920
$ws.Cells.Item(534,3).Value = 'Global'
$ws.Cells.Item(534,4).Value = 'M'
$ws.Cells.Item(534,5).Value = 'Sofiane Belle'
$ws.Cells.Item(534,6).Value = 'left forward'
$ws.Cells.Item(534,8).Value = 8.24
$ws.Cells.Item(534,9).Value = 1.32
$ws.Cells.Item(534,10).Value = 6.9
$ws.Cells.Item(534,11).Value = 0.73
$ws.Cells.Item(534,12).Value = 0.43
$ws.Cells.Item(534,13).Value = 0.15
$ws.Cells.Item(534,14).Value = 0.03
$ws.Cells.Item(534,15).Value = 9
$ws.Cells.Item(534,16).Value = 5.56
$ws.Cells.Item(534,17).Value = 31.07
$ws.Cells.Item(534,18).Value = 4.56
$ws.Cells.Item(534,19).Value = 21
$ws.Cells.Item(534,20).Value = 9
$ws.Cells.Item(534,21).Value = 36
$ws.Cells.Item(534,22).Value = 16

$ws.Cells.Item(535,2).Value = 45920
$ws.Cells.Item(535,3).Value = 'Global'
$ws.Cells.Item(535,4).Value = 'M'
$ws.Cells.Item(535,5).Value = 'Mattheo Haon'
$ws.Cells.Item(535,6).Value = 'right back'
$ws.Cells.Item(535,8).Value = 10.51
$ws.Cells.Item(535,9).Value = 1.7
$ws.Cells.Item(535,10).Value = 8.79
$ws.Cells.Item(535,11).Value = 1.03
$ws.Cells.Item(535,12).Value = 0.52
$ws.Cells.Item(535,13).Value = 0.15
$ws.Cells.Item(535,14).Value = 0.02
$ws.Cells.Item(535,15).Value = 12
$ws.Cells.Item(535,16).Value = 6.21
$ws.Cells.Item(535,17).Value = 32.19
$ws.Cells.Item(535,18).Value = 4.5
$ws.Cells.Item(535,19).Value = 24
$ws.Cells.Item(535,20).Value = 8
$ws.Cells.Item(535,21).Value = 27
$ws.Cells.Item(535,22).Value = 15

$ws.Cells.Item(536,2).Value = 45920
$ws.Cells.Item(536,3).Value = 'Global'
$ws.Cells.Item(536,4).Value = 'M'
$ws.Cells.Item(536,5).Value = 'Naim Ighbane'
$ws.Cells.Item(536,6).Value = 'center back'
$ws.Cells.Item(536,8).Value = 9.53
$ws.Cells.Item(536,9).Value = 1.19
$ws.Cells.Item(536,10).Value = 8.33
$ws.Cells.Item(536,11).Value = 0.73
$ws.Cells.Item(536,12).Value = 0.31
$ws.Cells.Item(536,13).Value = 0.15
$ws.Cells.Item(536,14).Value = 0.01
$ws.Cells.Item(536,15).Value = 8
$ws.Cells.Item(536,16).Value = 5.62
$ws.Cells.Item(536,17).Value = 30.63
$ws.Cells.Item(536,18).Value = 4.52
$ws.Cells.Item(536,19).Value = 34
$ws.Cells.Item(536,20).Value = 1
$ws.Cells.Item(536,21).Value = 27
$ws.Cells.Item(536,22).Value = 7

$ws.Cells.Item(537,2).Value = 45920
$ws.Cells.Item(537,3).Value = 'Global'
$ws.Cells.Item(537,4).Value = 'M'
$ws.Cells.Item(537,5).Value = 'Kamal Bafounta'
$ws.Cells.Item(537,6).Value = 'center midfield'
$ws.Cells.Item(537,8).Value = 10
$ws.Cells.Item(537,9).Value = 1.31
$ws.Cells.Item(537,10).Value = 8.67
$ws.Cells.Item(537,11).Value = 0.99
$ws.Cells.Item(537,12).Value = 0.22
$ws.Cells.Item(537,13).Value = 0.12
$ws.Cells.Item(537,14).Value = 0.01
$ws.Cells.Item(537,15).Value = 4
$ws.Cells.Item(537,16).Value = 6.02
$ws.Cells.Item(537,17).Value = 30.47
$ws.Cells.Item(537,18).Value = 4.28
$ws.Cells.Item(537,19).Value = 25
$ws.Cells.Item(537,20).Value = 1
$ws.Cells.Item(537,21).Value = 36
$ws.Cells.Item(537,22).Value = 7

$ws.Cells.Item(538,2).Value = 45920
$ws.Cells.Item(538,3).Value = 'Global'
$ws.Cells.Item(538,4).Value = 'M'
$ws.Cells.Item(538,5).Value = 'Karim Belmahi'
$ws.Cells.Item(538,6).Value = 'left forward'
$ws.Cells.Item(538,8).Value = 2.21
$ws.Cells.Item(538,9).Value = 0.51
$ws.Cells.Item(538,10).Value = 1.69
$ws.Cells.Item(538,11).Value = 0.37
$ws.Cells.Item(538,12).Value = 0.15
$ws.Cells.Item(538,13).Value = 0
$ws.Cells.Item(538,14).Value = 0
$ws.Cells.Item(538,15).Value = 0
$ws.Cells.Item(538,16).Value = 6.93
$ws.Cells.Item(538,17).Value = 24.88
$ws.Cells.Item(538,18).Value = 4.37
$ws.Cells.Item(538,19).Value = 14
$ws.Cells.Item(538,20).Value = 1
$ws.Cells.Item(538,21).Value = 14
$ws.Cells.Item(538,22).Value = 2

$ws.Cells.Item(539,2).Value = 45920
$ws.Cells.Item(539,3).Value = 'Global'
$ws.Cells.Item(539,4).Value = 'M'
$ws.Cells.Item(539,5).Value = 'Omar Benyounes'
$ws.Cells.Item(539,6).Value = 'center midfield'
$ws.Cells.Item(539,8).Value = 2.4
$ws.Cells.Item(539,9).Value = 0.42
$ws.Cells.Item(539,10).Value = 1.97
$ws.Cells.Item(539,11).Value = 0.29
$ws.Cells.Item(539,12).Value = 0.11
$ws.Cells.Item(539,13).Value = 0.03
$ws.Cells.Item(539,14).Value = 0
$ws.Cells.Item(539,15).Value = 2
$ws.Cells.Item(539,16).Value = 6.61
$ws.Cells.Item(539,17).Value = 29.12
$ws.Cells.Item(539,18).Value = 4.62
$ws.Cells.Item(539,19).Value = 3
$ws.Cells.Item(539,20).Value = 2
$ws.Cells.Item(539,21).Value = 14
$ws.Cells.Item(539,22).Value = 3

$ws.Cells.Item(540,2).Value = 45920
$ws.Cells.Item(540,3).Value = 'Global'
$ws.Cells.Item(540,4).Value = 'M'
$ws.Cells.Item(540,5).Value = 'Levy Ndoutoume'
$ws.Cells.Item(540,6).Value = 'left back'
$ws.Cells.Item(540,8).Value = 7.06
$ws.Cells.Item(540,9).Value = 1.25
$ws.Cells.Item(540,10).Value = 5.78
$ws.Cells.Item(540,11).Value = 0.83
$ws.Cells.Item(540,12).Value = 0.37
$ws.Cells.Item(540,13).Value = 0.07
$ws.Cells.Item(540,14).Value = 0
$ws.Cells.Item(540,15).Value = 6
$ws.Cells.Item(540,16).Value = 5.84
$ws.Cells.Item(540,17).Value = 28.84
$ws.Cells.Item(540,18).Value = 4.42
$ws.Cells.Item(540,19).Value = 30
$ws.Cells.Item(540,20).Value = 4
$ws.Cells.Item(540,21).Value = 39
$ws.Cells.Item(540,22).Value = 11

$ws.Cells.Item(541,2).Value = 45920
$ws.Cells.Item(541,3).Value = 'Global'
$ws.Cells.Item(541,4).Value = 'M'
$ws.Cells.Item(541,5).Value = 'Ilan Ihaddadene'
$ws.Cells.Item(541,6).Value = 'center midfield'
$ws.Cells.Item(541,8).Value = 9.77
$ws.Cells.Item(541,9).Value = 1.87
$ws.Cells.Item(541,10).Value = 7.87
$ws.Cells.Item(541,11).Value = 1.46
$ws.Cells.Item(541,12).Value = 0.35
$ws.Cells.Item(541,13).Value = 0.08
$ws.Cells.Item(541,14).Value = 0
$ws.Cells.Item(541,15).Value = 4
$ws.Cells.Item(541,16).Value = 7.17
$ws.Cells.Item(541,17).Value = 28.86
$ws.Cells.Item(541,18).Value = 4.85
$ws.Cells.Item(541,19).Value = 30
$ws.Cells.Item(541,20).Value = 5
$ws.Cells.Item(541,21).Value = 28
$ws.Cells.Item(541,22).Value = 6

$ws.Cells.Item(542,2).Value = 45922
$ws.Cells.Item(542,3).Value = 'Global'
$ws.Cells.Item(542,4).Value = 'J+2'
$ws.Cells.Item(542,5).Value = 'Fareh Wael'
$ws.Cells.Item(542,6).Value = 'center midfield'
$ws.Cells.Item(542,8).Value = 7.26
$ws.Cells.Item(542,9).Value = 0.02
$ws.Cells.Item(542,10).Value = 7.24
$ws.Cells.Item(542,11).Value = 0.02
$ws.Cells.Item(542,12).Value = 0
$ws.Cells.Item(542,13).Value = 0
$ws.Cells.Item(542,14).Value = 0
$ws.Cells.Item(542,15).Value = 0
$ws.Cells.Item(542,16).Value = 4.42
$ws.Cells.Item(542,17).Value = 17.04
$ws.Cells.Item(542,18).Value = 4.1
$ws.Cells.Item(542,19).Value = 11
$ws.Cells.Item(542,20).Value = 1
$ws.Cells.Item(542,21).Value = 4
$ws.Cells.Item(542,22).Value = 0

$ws.Cells.Item(543,2).Value = 45922
$ws.Cells.Item(543,3).Value = 'Global'
$ws.Cells.Item(543,4).Value = 'J+2'
$ws.Cells.Item(543,5).Value = 'Ilan Ihaddadene'
$ws.Cells.Item(543,6).Value = 'center midfield'
$ws.Cells.Item(543,8).Value = 8.07
$ws.Cells.Item(543,9).Value = 0.03
$ws.Cells.Item(543,10).Value = 8.05
$ws.Cells.Item(543,11).Value = 0.03
$ws.Cells.Item(543,12).Value = 0
$ws.Cells.Item(543,13).Value = 0
$ws.Cells.Item(543,14).Value = 0
$ws.Cells.Item(543,15).Value = 0
$ws.Cells.Item(543,16).Value = 5.1
$ws.Cells.Item(543,17).Value = 18.87
$ws.Cells.Item(543,18).Value = 4.37
$ws.Cells.Item(543,19).Value = 12
$ws.Cells.Item(543,20).Value = 1
$ws.Cells.Item(543,21).Value = 4
$ws.Cells.Item(543,22).Value = 0

$ws.Cells.Item(544,2).Value = 45922
$ws.Cells.Item(544,3).Value = 'Global'
$ws.Cells.Item(544,4).Value = 'J+2'
$ws.Cells.Item(544,5).Value = 'Emmanuel Valey'
$ws.Cells.Item(544,6).Value = 'left forward'
$ws.Cells.Item(544,8).Value = 8.72
$ws.Cells.Item(544,9).Value = 0.08
$ws.Cells.Item(544,10).Value = 8.63
$ws.Cells.Item(544,11).Value = 0.09
$ws.Cells.Item(544,12).Value = 0
$ws.Cells.Item(544,13).Value = 0
$ws.Cells.Item(544,14).Value = 0
$ws.Cells.Item(544,15).Value = 0
$ws.Cells.Item(544,16).Value = 4.29
$ws.Cells.Item(544,17).Value = 18.04
$ws.Cells.Item(544,18).Value = 5.01
$ws.Cells.Item(544,19).Value = 30
$ws.Cells.Item(544,20).Value = 5
$ws.Cells.Item(544,21).Value = 27
$ws.Cells.Item(544,22).Value = 8

$ws.Cells.Item(545,2).Value = 45922
$ws.Cells.Item(545,3).Value = 'Global'
$ws.Cells.Item(545,4).Value = 'J+2'
$ws.Cells.Item(545,5).Value = 'Hedi Nasri'
$ws.Cells.Item(545,6).Value = 'right back'
$ws.Cells.Item(545,8).Value = 7.44
$ws.Cells.Item(545,9).Value = 0.05
$ws.Cells.Item(545,10).Value = 7.39
$ws.Cells.Item(545,11).Value = 0.05
$ws.Cells.Item(545,12).Value = 0
$ws.Cells.Item(545,13).Value = 0
$ws.Cells.Item(545,14).Value = 0
$ws.Cells.Item(545,15).Value = 0
$ws.Cells.Item(545,16).Value = 4.66
$ws.Cells.Item(545,17).Value = 20.39
$ws.Cells.Item(545,18).Value = 3.95
$ws.Cells.Item(545,19).Value = 20
$ws.Cells.Item(545,20).Value = 0
$ws.Cells.Item(545,21).Value = 23
$ws.Cells.Item(545,22).Value = 2

$ws.Cells.Item(546,2).Value = 45922
$ws.Cells.Item(546,3).Value = 'Global'
$ws.Cells.Item(546,4).Value = 'J+2'
$ws.Cells.Item(546,5).Value = 'Mattheo Haon'
$ws.Cells.Item(546,6).Value = 'right back'
$ws.Cells.Item(546,8).Value = 7.43
$ws.Cells.Item(546,9).Value = 0.15
$ws.Cells.Item(546,10).Value = 7.28
$ws.Cells.Item(546,11).Value = 0.13
$ws.Cells.Item(546,12).Value = 0.03
$ws.Cells.Item(546,13).Value = 0
$ws.Cells.Item(546,14).Value = 0
$ws.Cells.Item(546,15).Value = 0
$ws.Cells.Item(546,16).Value = 4.7
$ws.Cells.Item(546,17).Value = 22.55
$ws.Cells.Item(546,18).Value = 4.8
$ws.Cells.Item(546,19).Value = 27
$ws.Cells.Item(546,20).Value = 2
$ws.Cells.Item(546,21).Value = 13
$ws.Cells.Item(546,22).Value = 2

$ws.Cells.Item(547,2).Value = 45922
$ws.Cells.Item(547,3).Value = 'Global'
$ws.Cells.Item(547,4).Value = 'J+2'
$ws.Cells.Item(547,5).Value = 'Malik Boussaid'
$ws.Cells.Item(547,6).Value = 'right back'
$ws.Cells.Item(547,8).Value = 7.87
$ws.Cells.Item(547,9).Value = 0.07
$ws.Cells.Item(547,10).Value = 7.8
$ws.Cells.Item(547,11).Value = 0.07
$ws.Cells.Item(547,12).Value = 0.01
$ws.Cells.Item(547,13).Value = 0
$ws.Cells.Item(547,14).Value = 0
$ws.Cells.Item(547,15).Value = 0
$ws.Cells.Item(547,16).Value = 4.57
$ws.Cells.Item(547,17).Value = 22.24
$ws.Cells.Item(547,18).Value = 4.31
$ws.Cells.Item(547,19).Value = 29
$ws.Cells.Item(547,20).Value = 3
$ws.Cells.Item(547,21).Value = 18
$ws.Cells.Item(547,22).Value = 4

$ws.Cells.Item(548,2).Value = 45923
$ws.Cells.Item(548,3).Value = 'Global'
$ws.Cells.Item(548,4).Value = 'J+3'
$ws.Cells.Item(548,5).Value = 'Kamal Bafounta'
$ws.Cells.Item(548,6).Value = 'center midfield'
$ws.Cells.Item(548,8).Value = 6.17
$ws.Cells.Item(548,9).Value = 0.22
$ws.Cells.Item(548,10).Value = 5.94
$ws.Cells.Item(548,11).Value = 0.21
$ws.Cells.Item(548,12).Value = 0.02
$ws.Cells.Item(548,13).Value = 0
$ws.Cells.Item(548,14).Value = 0
$ws.Cells.Item(548,15).Value = 0
$ws.Cells.Item(548,16).Value = 4.23
$ws.Cells.Item(548,17).Value = 23.05
$ws.Cells.Item(548,18).Value = 4.5
$ws.Cells.Item(548,19).Value = 27
$ws.Cells.Item(548,20).Value = 3
$ws.Cells.Item(548,21).Value = 23
$ws.Cells.Item(548,22).Value = 2

$ws.Cells.Item(549,2).Value = 45923
$ws.Cells.Item(549,3).Value = 'Global'
$ws.Cells.Item(549,4).Value = 'J+3'
$ws.Cells.Item(549,5).Value = 'Omar Benyounes'
$ws.Cells.Item(549,6).Value = 'center midfield'
$ws.Cells.Item(549,8).Value = 6.24
$ws.Cells.Item(549,9).Value = 0.19
$ws.Cells.Item(549,10).Value = 6.04
$ws.Cells.Item(549,11).Value = 0.17
$ws.Cells.Item(549,12).Value = 0.02
$ws.Cells.Item(549,13).Value = 0
$ws.Cells.Item(549,14).Value = 0
$ws.Cells.Item(549,15).Value = 0
$ws.Cells.Item(549,16).Value = 4.34
$ws.Cells.Item(549,17).Value = 22.51
$ws.Cells.Item(549,18).Value = 4.24
$ws.Cells.Item(549,19).Value = 31
$ws.Cells.Item(549,20).Value = 1
$ws.Cells.Item(549,21).Value = 19
$ws.Cells.Item(549,22).Value = 4

$ws.Cells.Item(550,2).Value = 45923
$ws.Cells.Item(550,3).Value = 'Global'
$ws.Cells.Item(550,4).Value = 'J+3'
$ws.Cells.Item(550,5).Value = 'Malik Boussaid'
$ws.Cells.Item(550,6).Value = 'right back'
$ws.Cells.Item(550,8).Value = 6.54
$ws.Cells.Item(550,9).Value = 0.25
$ws.Cells.Item(550,10).Value = 6.28
$ws.Cells.Item(550,11).Value = 0.25
$ws.Cells.Item(550,12).Value = 0.01
$ws.Cells.Item(550,13).Value = 0
$ws.Cells.Item(550,14).Value = 0
$ws.Cells.Item(550,15).Value = 0
$ws.Cells.Item(550,16).Value = 4.28
$ws.Cells.Item(550,17).Value = 21.93
$ws.Cells.Item(550,18).Value = 4.13
$ws.Cells.Item(550,19).Value = 47
$ws.Cells.Item(550,20).Value = 2
$ws.Cells.Item(550,21).Value = 32
$ws.Cells.Item(550,22).Value = 11

$ws.Cells.Item(551,2).Value = 45923
$ws.Cells.Item(551,3).Value = 'Global'
$ws.Cells.Item(551,4).Value = 'J+3'
$ws.Cells.Item(551,5).Value = 'Emmanuel Valey'
$ws.Cells.Item(551,6).Value = 'left forward'
$ws.Cells.Item(551,8).Value = 6.94
$ws.Cells.Item(551,9).Value = 0.16
$ws.Cells.Item(551,10).Value = 6.77
$ws.Cells.Item(551,11).Value = 0.16
$ws.Cells.Item(551,12).Value = 0.01
$ws.Cells.Item(551,13).Value = 0
$ws.Cells.Item(551,14).Value = 0
$ws.Cells.Item(551,15).Value = 0
$ws.Cells.Item(551,16).Value = 4.23
$ws.Cells.Item(551,17).Value = 21.28
$ws.Cells.Item(551,18).Value = 4.66
$ws.Cells.Item(551,19).Value = 37
$ws.Cells.Item(551,20).Value = 4
$ws.Cells.Item(551,21).Value = 37
$ws.Cells.Item(551,22).Value = 11

$ws.Cells.Item(552,2).Value = 45923
$ws.Cells.Item(552,3).Value = 'Global'
$ws.Cells.Item(552,4).Value = 'J+3'
$ws.Cells.Item(552,5).Value = 'Hedi Nasri'
$ws.Cells.Item(552,6).Value = 'right back'
$ws.Cells.Item(552,8).Value = 6.28
$ws.Cells.Item(552,9).Value = 0.2
$ws.Cells.Item(552,10).Value = 6.07
$ws.Cells.Item(552,11).Value = 0.18
$ws.Cells.Item(552,12).Value = 0.02
$ws.Cells.Item(552,13).Value = 0
$ws.Cells.Item(552,14).Value = 0
$ws.Cells.Item(552,15).Value = 0
$ws.Cells.Item(552,16).Value = 4.44
$ws.Cells.Item(552,17).Value = 23.32
$ws.Cells.Item(552,18).Value = 4.61
$ws.Cells.Item(552,19).Value = 22
$ws.Cells.Item(552,20).Value = 5
$ws.Cells.Item(552,21).Value = 16
$ws.Cells.Item(552,22).Value = 3

# Step 6: apply styles to B (date format) and D (centered) columns by copying from row 527
$ws.Range("B527").Copy()
$ws.Cells.Item(528,2).PasteSpecial(-4122)
$ws.Cells.Item(529,2).PasteSpecial(-4122)
$ws.Cells.Item(530,2).PasteSpecial(-4122)
$ws.Cells.Item(531,2).PasteSpecial(-4122)
$ws.Cells.Item(532,2).PasteSpecial(-4122)
$ws.Cells.Item(533,2).PasteSpecial(-4122)
$ws.Cells.Item(534,2).PasteSpecial(-4122)
$ws.Cells.Item(535,2).PasteSpecial(-4122)
$ws.Cells.Item(536,2).PasteSpecial(-4122)
$ws.Cells.Item(537,2).PasteSpecial(-4122)
$ws.Cells.Item(538,2).PasteSpecial(-4122)
$ws.Cells.Item(539,2).PasteSpecial(-4122)
$ws.Cells.Item(540,2).PasteSpecial(-4122)
$ws.Cells.Item(541,2).PasteSpecial(-4122)
$ws.Cells.Item(542,2).PasteSpecial(-4122)
$ws.Cells.Item(543,2).PasteSpecial(-4122)
$ws.Cells.Item(544,2).PasteSpecial(-4122)
$ws.Cells.Item(545,2).PasteSpecial(-4122)
$ws.Cells.Item(546,2).PasteSpecial(-4122)
$ws.Cells.Item(547,2).PasteSpecial(-4122)
$ws.Cells.Item(548,2).PasteSpecial(-4122)
$ws.Cells.Item(549,2).PasteSpecial(-4122)
$ws.Cells.Item(550,2).PasteSpecial(-4122)
$ws.Cells.Item(551,2).PasteSpecial(-4122)
$ws.Cells.Item(552,2).PasteSpecial(-4122)
$ws.Range("D527").Copy()
$ws.Cells.Item(528,4).PasteSpecial(-4122)
$ws.Cells.Item(529,4).PasteSpecial(-4122)
$ws.Cells.Item(530,4).PasteSpecial(-4122)
$ws.Cells.Item(531,4).PasteSpecial(-4122)
$ws.Cells.Item(532,4).PasteSpecial(-4122)
$ws.Cells.Item(533,4).PasteSpecial(-4122)
$ws.Cells.Item(534,4).PasteSpecial(-4122)
$ws.Cells.Item(535,4).PasteSpecial(-4122)
$ws.Cells.Item(536,4).PasteSpecial(-4122)
$ws.Cells.Item(537,4).PasteSpecial(-4122)
$ws.Cells.Item(538,4).PasteSpecial(-4122)
$ws.Cells.Item(539,4).PasteSpecial(-4122)
$ws.Cells.Item(540,4).PasteSpecial(-4122)
$ws.Cells.Item(541,4).PasteSpecial(-4122)
$ws.Cells.Item(542,4).PasteSpecial(-4122)
$ws.Cells.Item(543,4).PasteSpecial(-4122)
$ws.Cells.Item(544,4).PasteSpecial(-4122)
$ws.Cells.Item(545,4).PasteSpecial(-4122)
$ws.Cells.Item(546,4).PasteSpecial(-4122)
$ws.Cells.Item(547,4).PasteSpecial(-4122)
$ws.Cells.Item(548,4).PasteSpecial(-4122)
$ws.Cells.Item(549,4).PasteSpecial(-4122)
$ws.Cells.Item(550,4).PasteSpecial(-4122)
$ws.Cells.Item(551,4).PasteSpecial(-4122)
$ws.Cells.Item(552,4).PasteSpecial(-4122)

# Step 7: update view selection to match target
$ws.Range("B555").Select()
